# Updated symbol list on Wed Dec 21 17:51:18 UTC 2022 with GitHub Actions
#
# The "Price" column (D) in this sheet is stored as TEXT (inline strings
# that merely look like numbers), not as actual numbers. A plain
# `$range.Value = "123.45"` assignment lets Excel's normal type inference
# kick in and silently turns the cell into a real Number, which would
# change both the stored type and the cell's formatting/style record.
# To faithfully reproduce the source edit (text in, text out, no style
# churn) we momentarily force a "Text" number format before writing the
# literal string, then clear the formatting override again so the cell's
# style reverts to the same (unset/default) style it had before.

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates -------------------------------------------------
Set-TextValue $ws "D2"  "247.77"
Set-TextValue $ws "D3"  "22.43"
Set-TextValue $ws "D5"  "0.05687"
Set-TextValue $ws "D6"  "3.419"
Set-TextValue $ws "D7"  "6.310"
Set-TextValue $ws "D8"  "0.8067"
Set-TextValue $ws "D9"  "0.8996"
Set-TextValue $ws "D10" "0.1423"
Set-TextValue $ws "D11" "0.07429"
Set-TextValue $ws "D12" "0.03056"
Set-TextValue $ws "D13" "0.03076"
Set-TextValue $ws "D14" "0.09395"
Set-TextValue $ws "D15" "3.893"
Set-TextValue $ws "D16" "0.001594"
Set-TextValue $ws "D17" "0.04783"
Set-TextValue $ws "D18" "0.01827"
Set-TextValue $ws "D19" "0.0005812"
Set-TextValue $ws "D20" "0.006412"
Set-TextValue $ws "D21" "0.005040"
Set-TextValue $ws "D22" "0.0009970"
Set-TextValue $ws "D23" "0.0001500"
Set-TextValue $ws "D24" "3.695"
Set-TextValue $ws "D25" "2.160"
Set-TextValue $ws "D26" "0.3247"
Set-TextValue $ws "D27" "0.1341"
Set-TextValue $ws "D40" "0.03957"
Set-TextValue $ws "D41" "0.003043"
Set-TextValue $ws "D43" "0.002731"
Set-TextValue $ws "D44" "0.008759"
Set-TextValue $ws "D45" "0.00005580"
Set-TextValue $ws "D47" "0.4992"
Set-TextValue $ws "D48" "0.1951"

# --- Column E ("Volume(1h)") updates --------------------------------------------
# These are plain alphanumeric labels (not numeric-looking), so a normal
# Value assignment keeps them as text without touching formatting.
$ws.Range("E19").Value = "18OneONE"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
